$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values for column G, rows 2-40, replacing old Strike# values
$newValues = @{
    2 = 8
    3 = 7
    4 = 7
    5 = 9
    6 = 7
    7 = 5
    8 = 9
    9 = 6
    10 = 5
    11 = 4
    12 = 4
    13 = 4
    14 = 4
    15 = 6
    16 = 4
    17 = 9
    18 = 8
    19 = 9
    20 = 6
    21 = 5
    22 = 2
    23 = 4
    24 = 6
    25 = 5
    26 = 5
    27 = 8
    28 = 7
    29 = 4
    30 = 10
    31 = 7
    32 = 9
    33 = 3
    34 = 4
    35 = 3
    36 = 4
    37 = 5
    38 = 5
    39 = 0
    40 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
